$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1067073025158294
$ws.Range("D2").Value = 0.01745812549005166
$ws.Range("E2").Value = 0.4207214845885545
$ws.Range("F2").Value = 0.5547355846707873
$ws.Range("G2").Value = 0.002395864625174338
$ws.Range("I2").Value = 0.3969603734482963
$ws.Range("K2").Value = 0.8925442119360696
$ws.Range("O2").Value = 1.816893039332911
$ws.Range("B3").Value = 0.09439299549843838
$ws.Range("D3").Value = 0.01595942225532809
$ws.Range("E3").Value = 0.3669982134373555
$ws.Range("F3").Value = 0.5503303348038884
$ws.Range("G3").Value = 0.002398822723454029
$ws.Range("I3").Value = 0.4047372261755875
$ws.Range("K3").Value = 0.7786953463426869
$ws.Range("O3").Value = 1.816724314547827
$ws.Range("B4").Value = 0.08682080916184987
$ws.Range("D4").Value = 0.01503396172898164
$ws.Range("E4").Value = 0.3341047463456306
$ws.Range("F4").Value = 0.5481239848103101
$ws.Range("G4").Value = 0.002400733500092904
$ws.Range("I4").Value = 0.4097942213140899
$ws.Range("K4").Value = 0.7085436272149082
$ws.Range("O4").Value = 1.818285711413893
$ws.Range("B5").Value = 0.0837325212358877
$ws.Range("D5").Value = 0.01465553884172976
$ws.Range("E5").Value = 0.3207218923371187
$ws.Range("F5").Value = 0.5473497653400727
$ws.Range("G5").Value = 0.002401535994890106
$ws.Range("I5").Value = 0.4119258397448489
$ws.Range("K5").Value = 0.6798950155525176
$ws.Range("O5").Value = 1.819339002806061
$ws.Range("B6").Value = 0.08321956594035385
$ws.Range("D6").Value = 0.01459262497788671
$ws.Range("E6").Value = 0.3185009245696335
$ws.Range("F6").Value = 0.547228736054457
$ws.Range("G6").Value = 0.002401670690390636
$ws.Range("I6").Value = 0.4122840710907854
$ws.Range("K6").Value = 0.675134272577111
$ws.Range("O6").Value = 1.819539038681256
$ws.Range("B7").Value = 0.08677916947394237
$ws.Range("D7").Value = 0.01502886337021891
$ws.Range("E7").Value = 0.3339241758557989
$ws.Range("F7").Value = 0.5481130383741402
$ws.Range("G7").Value = 0.002400744226187579
$ws.Range("I7").Value = 0.4098226822603288
$ws.Range("K7").Value = 0.7081575083314533
$ws.Range("O7").Value = 1.818298230235769
$ws.Range("B8").Value = 0.1024637859169246
$ws.Range("D8").Value = 0.01694247989700415
$ws.Range("E8").Value = 0.4021773040717562
$ws.Range("F8").Value = 0.5531129102255079
$ws.Range("G8").Value = 0.002396865012055266
$ws.Range("I8").Value = 0.3995832767112169
$ws.Range("K8").Value = 0.8533412636975299
$ws.Range("O8").Value = 1.816488282883938
$ws.Range("B9").Value = 0.1331236249155694
$ws.Range("D9").Value = 0.02065227895133148
$ws.Range("E9").Value = 0.5368476431241476
$ws.Range("F9").Value = 0.5668953387405438
$ws.Range("G9").Value = 0.002390004114536529
$ws.Range("I9").Value = 0.3817433259562462
$ws.Range("K9").Value = 1.136041209528173
$ws.Range("O9").Value = 1.826228267593962
$ws.Range("B10").Value = 0.1555798490400093
$ws.Range("D10").Value = 0.02335052844927787
$ws.Range("E10").Value = 0.6364267113224997
$ws.Range("F10").Value = 0.5794793289163636
$ws.Range("G10").Value = 0.002385413383827893
$ws.Range("I10").Value = 0.3700041448694824
$ws.Range("K10").Value = 1.342490120229627
$ws.Range("O10").Value = 1.841598627115303
$ws.Range("B11").Value = 0.1657787032538351
$ws.Range("D11").Value = 0.02457185370379733
$ws.Range("E11").Value = 0.6818955294875337
$ws.Range("F11").Value = 0.5857451354050767
$ws.Range("G11").Value = 0.002383421591611978
$ws.Range("I11").Value = 0.3649612562917985
$ws.Range("K11").Value = 1.436132998613346
$ws.Range("O11").Value = 1.850399360790391
$ws.Range("B12").Value = 0.1696381471117263
$ws.Range("D12").Value = 0.0250334320915826
$ws.Range("E12").Value = 0.6991400252582594
$ws.Range("F12").Value = 0.5881962200375881
$ws.Range("G12").Value = 0.002382681156023435
$ws.Range("I12").Value = 0.3630944621362344
$ws.Range("K12").Value = 1.471553223464639
$ws.Range("O12").Value = 1.853993988466584
$ws.Range("B13").Value = 0.1688070683075296
$ws.Range("D13").Value = 0.02493406392626696
$ws.Range("E13").Value = 0.6954249061123363
$ws.Range("F13").Value = 0.5876648416832069
$ws.Range("G13").Value = 0.00238284000886157
$ws.Range("I13").Value = 0.3634946030737627
$ws.Range("K13").Value = 1.463926655938224
$ws.Range("O13").Value = 1.853208139150411
$ws.Range("B14").Value = 0.1660962763995713
$ws.Range("D14").Value = 0.02460984646693731
$ws.Range("E14").Value = 0.6833137029393583
$ws.Range("F14").Value = 0.5859452141637007
$ws.Range("G14").Value = 0.002383360399082846
$ws.Range("I14").Value = 0.3648068149879977
$ws.Range("K14").Value = 1.43904785644952
$ws.Range("O14").Value = 1.850689831605649
$ws.Range("B15").Value = 0.1644354866439386
$ws.Range("D15").Value = 0.02441113437340192
$ws.Range("E15").Value = 0.6758987451607368
$ws.Range("F15").Value = 0.5849021129645706
$ws.Range("G15").Value = 0.002383680949955766
$ws.Range("I15").Value = 0.365616163981028
$ws.Range("K15").Value = 1.423803589103841
$ws.Range("O15").Value = 1.849181469528219
$ws.Range("B16").Value = 0.1549129741910065
$ws.Range("D16").Value = 0.02327058615215805
$ws.Range("E16").Value = 0.6334588064632953
$ws.Range("F16").Value = 0.5790807804910827
$ws.Range("G16").Value = 0.002385545488844951
$ws.Range("I16").Value = 0.3703396995288468
$ws.Range("K16").Value = 1.3363647644602
$ws.Range("O16").Value = 1.841060035535122
$ws.Range("B17").Value = 0.1490667890881667
$ws.Range("D17").Value = 0.02256930682602842
$ws.Range("E17").Value = 0.6074682018785893
$ws.Range("F17").Value = 0.5756485898308767
$ws.Range("G17").Value = 0.00238671400178579
$ws.Range("I17").Value = 0.3733136464589668
$ws.Range("K17").Value = 1.282653397138802
$ws.Range("O17").Value = 1.83654239095722
$ws.Range("B18").Value = 0.14570267114442
$ws.Range("D18").Value = 0.02216537507850802
$ws.Range("E18").Value = 0.5925349537128994
$ws.Range("F18").Value = 0.5737254015578657
$ws.Range("G18").Value = 0.002387395192158127
$ws.Range("I18").Value = 0.3750521625396726
$ws.Range("K18").Value = 1.251734531184297
$ws.Range("O18").Value = 1.834114089003009
$ws.Range("B19").Value = 0.1445633807234543
$ws.Range("D19").Value = 0.02202851301286302
$ws.Range("E19").Value = 0.5874814751355615
$ws.Range("F19").Value = 0.573082972185432
$ws.Range("G19").Value = 0.002387627395175848
$ws.Range("I19").Value = 0.3756455973574582
$ws.Range("K19").Value = 1.241261590894283
$ws.Range("O19").Value = 1.833321069696723
$ws.Range("B20").Value = 0.1496892874222624
$ws.Range("D20").Value = 0.02264401888158574
$ws.Range("E20").Value = 0.610233290068507
$ws.Range("F20").Value = 0.5760086788764056
$ws.Range("G20").Value = 0.002386588671169674
$ws.Range("I20").Value = 0.3729941682551816
$ws.Range("K20").Value = 1.288373712781436
$ws.Range("O20").Value = 1.837005680638129
$ws.Range("B21").Value = 0.1668925753849777
$ws.Range("D21").Value = 0.02470510195967535
$ws.Range("E21").Value = 0.6868703230467332
$ws.Range("F21").Value = 0.5864481796104997
$ws.Range("G21").Value = 0.002383207173545797
$ws.Range("I21").Value = 0.3644202232507476
$ws.Range("K21").Value = 1.446356466327245
$ws.Range("O21").Value = 1.851422394349811
$ws.Range("B22").Value = 0.1781204088217407
$ws.Range("D22").Value = 0.02604681510985785
$ws.Range("E22").Value = 0.7371121287652045
$ws.Range("F22").Value = 0.5937279813233545
$ws.Range("G22").Value = 0.002381077648015092
$ws.Range("I22").Value = 0.3590663883962328
$ws.Range("K22").Value = 1.549371870185553
$ws.Range("O22").Value = 1.862372460073715
$ws.Range("B23").Value = 0.1721294055082865
$ws.Range("D23").Value = 0.02533121468712807
$ws.Range("E23").Value = 0.7102822675780374
$ws.Range("F23").Value = 0.5898006266316429
$ws.Range("G23").Value = 0.002382206875666704
$ws.Range("I23").Value = 0.3619009523161001
$ws.Range("K23").Value = 1.494412553646328
$ws.Range("O23").Value = 1.856387771911159
$ws.Range("B24").Value = 0.1494078654641697
$ws.Range("D24").Value = 0.02261024391665245
$ws.Range("E24").Value = 0.6089831655574329
$ws.Range("F24").Value = 0.575845726831858
$ws.Range("G24").Value = 0.002386645303996115
$ws.Range("I24").Value = 0.3731385147876134
$ws.Range("K24").Value = 1.285787680990097
$ws.Range("O24").Value = 1.83679570121646
$ws.Range("B25").Value = 0.124840825348457
$ws.Range("D25").Value = 0.01965339264135935
$ws.Range("E25").Value = 0.5003136139544608
$ws.Range("F25").Value = 0.5627375683187523
$ws.Range("G25").Value = 0.002391780793488951
$ws.Range("I25").Value = 0.3863294898702501
$ws.Range("K25").Value = 1.059781038383164
$ws.Range("O25").Value = 1.822158992178856
